$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 2378.111
$ws.Range("I2").Value = 2612.875
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 2612.875
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -2499.875
$ws.Range("N2").Value = -726
# Row 6
$ws.Range("H6").Value = 14445.167
$ws.Range("I6").Value = 14445.167
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 43335.501
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -43223.501
$ws.Range("N6").ClearContents()
# Row 10
$ws.Range("H10").Value = 4624.5
$ws.Range("J10").Value = 4624.5
$ws.Range("L10").Value = 4624.5
$ws.Range("N10").Value = -5210.5
# Row 19
$ws.Range("H19").Value = 1418.25
$ws.Range("I19").Value = 741.8
$ws.Range("J19").Value = 1725.7273
$ws.Range("K19").Value = 741.8
$ws.Range("L19").Value = 1725.7273
$ws.Range("M19").Value = -566.8
$ws.Range("N19").Value = -2075.7273
# Row 29
$ws.Range("H29").Value = 4531.875
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 5036.4287
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 15109.2861
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = -15671.2861
# Row 62
$ws.Range("H62").Value = 8425.444
$ws.Range("I62").Value = 1715.9
$ws.Range("J62").Value = 16812.375
$ws.Range("K62").Value = 1715.9
$ws.Range("L62").Value = 16812.375
$ws.Range("M62").Value = -1091.9
$ws.Range("N62").Value = -18060.375
# Row 65
$ws.Range("H65").Value = 8425.444
$ws.Range("I65").Value = 1715.9
$ws.Range("J65").Value = 16812.375
$ws.Range("K65").Value = 8579.5
$ws.Range("L65").Value = 84061.875
$ws.Range("M65").Value = -5459.5
$ws.Range("N65").Value = -90301.875
# Row 76
$ws.Range("H76").Value = 6012.8184
$ws.Range("I76").Value = 4228.2
$ws.Range("K76").Value = 4228.2
$ws.Range("M76").Value = -3913.2
# Row 79
$ws.Range("H79").Value = 6012.8184
$ws.Range("I79").Value = 4228.2
$ws.Range("K79").Value = 4228.2
$ws.Range("M79").Value = -3136.2
# Row 116
$ws.Range("H116").Value = 6454.524
$ws.Range("I116").Value = 6076.6665
$ws.Range("K116").Value = 6076.6665
$ws.Range("M116").Value = -2634.6665
# Row 132
$ws.Range("H132").Value = 1570.2941
$ws.Range("I132").Value = 1570.2941
$ws.Range("K132").Value = 4710.8823
$ws.Range("M132").Value = -2180.8823
# Row 135
$ws.Range("H135").Value = 2966.8333
$ws.Range("I135").Value = 2911.2
$ws.Range("J135").Value = 3245
$ws.Range("K135").Value = 26200.8
$ws.Range("L135").Value = 29205
$ws.Range("M135").Value = -23665.8
$ws.Range("N135").Value = -34275
# Row 138
$ws.Range("H138").Value = 1884.1818
$ws.Range("I138").Value = 1633.2727
$ws.Range("J138").Value = 1915.5454
$ws.Range("K138").Value = 4899.8181
$ws.Range("L138").Value = 5746.6362
$ws.Range("M138").Value = 240.1818999999996
$ws.Range("N138").Value = -16026.6362

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 540.6539
$ws.Range("I2").Value = 420.81818
$ws.Range("J2").Value = 1199.75
$ws.Range("K2").Value = 420.81818
$ws.Range("L2").Value = 1199.75
$ws.Range("M2").Value = -307.81818
$ws.Range("N2").Value = -1425.75
# Row 32
$ws.Range("H32").Value = 9624174
$ws.Range("I32").Value = 12201419
$ws.Range("J32").Value = 18076.908
$ws.Range("K32").Value = 12201419
$ws.Range("L32").Value = 18076.908
$ws.Range("M32").Value = -12201132
$ws.Range("N32").Value = -18650.908
# Row 61
$ws.Range("H61").Value = 22732114
$ws.Range("I61").Value = 16670917
$ws.Range("K61").Value = 16670917
$ws.Range("M61").Value = -16670705
# Row 110
$ws.Range("H110").Value = 1731.2
$ws.Range("I110").Value = 1573.4615
$ws.Range("K110").Value = 1573.4615
$ws.Range("M110").Value = 471.5385000000001
# Row 116
$ws.Range("H116").Value = 540.6539
$ws.Range("I116").Value = 420.81818
$ws.Range("J116").Value = 1199.75
$ws.Range("K116").Value = 420.81818
$ws.Range("L116").Value = 1199.75
$ws.Range("M116").Value = 1873.18182
$ws.Range("N116").Value = -5787.75
# Row 132
$ws.Range("H132").Value = 8526
$ws.Range("I132").Value = 3390
$ws.Range("K132").Value = 10170
$ws.Range("M132").Value = -7640
# Row 136
$ws.Range("H136").Value = 22732114
$ws.Range("I136").Value = 16670917
$ws.Range("K136").Value = 50012751
$ws.Range("M136").Value = -50010201

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 540.6539
$ws.Range("I3").Value = 420.81818
$ws.Range("J3").Value = 1199.75
$ws.Range("K3").Value = 420.81818
$ws.Range("L3").Value = 1199.75
$ws.Range("M3").Value = -306.81818
$ws.Range("N3").Value = -1427.75
# Row 105
$ws.Range("H105").Value = 2904.2632
$ws.Range("I105").Value = 4016.3333
$ws.Range("K105").Value = 4016.3333
$ws.Range("M105").Value = -2269.3333
# Row 107
$ws.Range("H107").Value = 1060.9
$ws.Range("I107").Value = 970.94116
$ws.Range("J107").Value = 1570.6666
$ws.Range("K107").Value = 970.94116
$ws.Range("L107").Value = 1570.6666
$ws.Range("M107").Value = 949.05884
$ws.Range("N107").Value = -5410.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 455105.25
$ws.Range("I31").Value = 6693.5137
$ws.Range("K31").Value = 6693.5137
$ws.Range("M31").Value = -6398.5137
# Row 34
$ws.Range("H34").Value = 455105.25
$ws.Range("I34").Value = 6693.5137
$ws.Range("K34").Value = 6693.5137
$ws.Range("M34").Value = -6491.5137
# Row 69
$ws.Range("H69").Value = 103664
$ws.Range("I69").Value = 99496
$ws.Range("J69").Value = 112000
$ws.Range("K69").Value = 99496
$ws.Range("L69").Value = 112000
$ws.Range("M69").Value = -98747
$ws.Range("N69").Value = -113498
# Row 72
$ws.Range("H72").Value = 103664
$ws.Range("I72").Value = 99496
$ws.Range("J72").Value = 112000
$ws.Range("K72").Value = 298488
$ws.Range("L72").Value = 336000
$ws.Range("M72").Value = -294744
$ws.Range("N72").Value = -343488
# Row 122
$ws.Range("H122").Value = 2490.2856
$ws.Range("I122").Value = 2767
$ws.Range("K122").Value = 8301
$ws.Range("M122").Value = -5851

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 20000
$ws.Range("J3").Value = 20000
$ws.Range("L3").Value = 60000
$ws.Range("N3").Value = -60224
# Row 5
$ws.Range("H5").Value = 1975.2
$ws.Range("I5").Value = 1975.2
$ws.Range("K5").Value = 5925.6
$ws.Range("M5").Value = -5813.6
# Row 119
$ws.Range("H119").Value = 9444.223
$ws.Range("I119").Value = 2499.5
$ws.Range("K119").Value = 7498.5
$ws.Range("M119").Value = -2660.5
# Row 125
$ws.Range("H125").Value = 4954.5454
$ws.Range("I125").Value = 2500
$ws.Range("J125").Value = 5200
$ws.Range("K125").Value = 7500
$ws.Range("L125").Value = 15600
$ws.Range("M125").Value = -2580
$ws.Range("N125").Value = -25440
# Row 135
$ws.Range("H135").Value = 1975.2
$ws.Range("I135").Value = 1975.2
$ws.Range("K135").Value = 17776.8
$ws.Range("M135").Value = -15241.8

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 12609.375
$ws.Range("I13").Value = 9410.714
$ws.Range("J13").Value = 35000
$ws.Range("K13").Value = 9410.714
$ws.Range("L13").Value = 35000
$ws.Range("M13").Value = -9271.714
$ws.Range("N13").Value = -35278
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 97
$ws.Range("H97").Value = 848.375
$ws.Range("J97").Value = 791.375
$ws.Range("L97").Value = 791.375
$ws.Range("N97").Value = -1783.375

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1346.875
$ws.Range("I22").Value = 1212.1666
$ws.Range("K22").Value = 1212.1666
$ws.Range("M22").Value = -917.1666
# Row 27
$ws.Range("H27").Value = 1346.875
$ws.Range("I27").Value = 1212.1666
$ws.Range("K27").Value = 1212.1666
$ws.Range("M27").Value = -1105.1666
# Row 100
$ws.Range("H100").Value = 3894.4
$ws.Range("I100").Value = 3250
$ws.Range("J100").Value = 4055.5
$ws.Range("K100").Value = 3250
$ws.Range("L100").Value = 4055.5
$ws.Range("M100").Value = -2709
$ws.Range("N100").Value = -5137.5
# Row 109
$ws.Range("H109").Value = 99425
$ws.Range("J109").Value = 99425
$ws.Range("L109").Value = 99425
$ws.Range("N109").Value = -102199
# Row 122
$ws.Range("H122").Value = 5370.5586
$ws.Range("I122").Value = 4749.087
$ws.Range("K122").Value = 14247.261
$ws.Range("M122").Value = -11797.261
# Row 123
$ws.Range("H123").Value = 53330
$ws.Range("J123").Value = 53330
$ws.Range("L123").Value = 53330
$ws.Range("N123").Value = -63130
# Row 132
$ws.Range("H132").Value = 849817.25
$ws.Range("I132").Value = 23394.125
$ws.Range("K132").Value = 70182.375
$ws.Range("M132").Value = -67652.375

$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 266648.28
$ws.Range("I31").Value = 773750
$ws.Range("J31").Value = 63807.6
$ws.Range("K31").Value = 773750
$ws.Range("L31").Value = 63807.6
$ws.Range("M31").Value = -773402
$ws.Range("N31").Value = -64503.6
# Row 109
$ws.Range("H109").Value = 104980
$ws.Range("J109").Value = 104980
$ws.Range("L109").Value = 104980
$ws.Range("N109").Value = -107754
# Row 122
$ws.Range("H122").Value = 2238.2917
$ws.Range("I122").Value = 2032.25
$ws.Range("K122").Value = 6096.75
$ws.Range("M122").Value = -3646.75
